# Adds a second "Note 2" paragraph (with a new empty paragraph before it)
# right after the existing "Note : ... main." paragraph, matching the
# commit "More theory on how to work with Git".
#
# The insertion point is anchored on the document's existing "_GoBack"
# bookmark, which in the source document sits right at the end of the
# "Note :" paragraph (immediately after "main." and before the paragraph
# mark). That is also where the new content - and a freshly relocated
# "_GoBack" bookmark - ends up after the edit, per the target diff.

$d = $word.ActiveDocument

$bm = $d.Bookmarks.Item("_GoBack")
$insertPos = $bm.Start
$bm.Delete()

$r = $d.Range($insertPos, $insertPos)

$xmlFrag = '<w:p/><w:p><w:pPr><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve">Note </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">2 </w:t></w:r><w:r><w:t>:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> In official projects, The </w:t></w:r><w:r><w:rPr><w:color w:val="833C0B" w:themeColor="accent2" w:themeShade="80"/></w:rPr><w:t>main</w:t></w:r><w:r><w:rPr><w:color w:val="833C0B" w:themeColor="accent2" w:themeShade="80"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>branch must be empty until the project is assuredly over.</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> We instead work on the project in using</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> other branches.</w:t></w:r></w:p>'

$r.InsertXML($xmlFrag)
